$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update task estimates (Week 1 column) for rows affected by NPC refactor/implementation progress
$ws.Range("C3").Value = 1.5
$ws.Range("C4").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0

# Update the selected cell in the sheet view
$ws.Range("B20").Select()

$wb.Save()
